$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Updated match-by-match stats (runs, balls, fours, sixes) for Priyam Garg
# "updated activity till excel form" - rows refreshed with latest data.
Set-TextValue "C2" "7"
Set-TextValue "D2" "14"
Set-TextValue "E2" "0"
Set-TextValue "F2" "0"

Set-TextValue "C3" "3"
Set-TextValue "D3" "5"
Set-TextValue "E3" "0"
Set-TextValue "F3" "0"

Set-TextValue "C4" "17"
Set-TextValue "D4" "12"
Set-TextValue "E4" "0"
Set-TextValue "F4" "2"

Set-TextValue "C5" "16"
Set-TextValue "D5" "18"
Set-TextValue "E5" "1"
Set-TextValue "F5" "0"

Set-TextValue "C6" "4"
Set-TextValue "D6" "7"
Set-TextValue "E6" "0"
Set-TextValue "F6" "0"

Set-TextValue "C7" "15"
Set-TextValue "D7" "8"
Set-TextValue "E7" "1"
Set-TextValue "F7" "1"

Set-TextValue "C10" "0"
Set-TextValue "D10" "1"
Set-TextValue "E10" "0"
Set-TextValue "F10" "0"

Set-TextValue "C11" "12"
Set-TextValue "D11" "13"
Set-TextValue "E11" "1"
Set-TextValue "F11" "0"
